# Update "想去人数" (want-to-go count) values in column F across the
# four worksheets of the workbook, per the commit's regenerated data.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (sheet1) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1367
$ws1.Range("F9").Value  = 57
$ws1.Range("F10").Value = 8592
$ws1.Range("F11").Value = 483
$ws1.Range("F13").Value = 133
$ws1.Range("F14").Value = 112
$ws1.Range("F15").Value = 291
$ws1.Range("F16").Value = 331
$ws1.Range("F17").Value = 88
$ws1.Range("F19").Value = 10861
$ws1.Range("F20").Value = 289
$ws1.Range("F24").Value = 10
$ws1.Range("F26").Value = 91
$ws1.Range("F28").Value = 2678
$ws1.Range("F30").Value = 37
$ws1.Range("F33").Value = 891
$ws1.Range("F34").Value = 4078
$ws1.Range("F35").Value = 2539
$ws1.Range("F36").Value = 279
$ws1.Range("F37").Value = 2580
$ws1.Range("F39").Value = 1243
$ws1.Range("F41").Value = 754
$ws1.Range("F43").Value = 316
$ws1.Range("F49").Value = 83

# ---- Sheet "演出" (sheet2) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 14
$ws2.Range("F5").Value  = 6
$ws2.Range("F7").Value  = 45
$ws2.Range("F21").Value = 30

# ---- Sheet "本地生活" (sheet3) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 33

# ---- Sheet "全部类型" (sheet4) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6
$ws4.Range("F6").Value  = 1367
$ws4.Range("F11").Value = 8592
$ws4.Range("F12").Value = 483
$ws4.Range("F13").Value = 112
$ws4.Range("F14").Value = 291
$ws4.Range("F15").Value = 331
$ws4.Range("F16").Value = 88
$ws4.Range("F18").Value = 10862
$ws4.Range("F19").Value = 289
$ws4.Range("F20").Value = 33
$ws4.Range("F24").Value = 10
$ws4.Range("F27").Value = 91
$ws4.Range("F29").Value = 2678
$ws4.Range("F32").Value = 891
$ws4.Range("F34").Value = 4078
$ws4.Range("F35").Value = 2539
$ws4.Range("F36").Value = 279
$ws4.Range("F37").Value = 2580
$ws4.Range("F40").Value = 1243
$ws4.Range("F42").Value = 754
$ws4.Range("F44").Value = 316
$ws4.Range("F49").Value = 83
